# Applies the 05-11-2023 08:45 scraper refresh to the Switzerland Super League 2023-2024 sheet:
# - several already-recorded fixtures had been associated with the wrong match in the pair
#   played at the same kickoff time; columns F:V (teams, scores, odds, timestamps, match url)
#   are corrected in place for those rows (A:E - index/country/tourney/season/date - are untouched).
# - three newly played fixtures (matchweek of 29/10-04/11/2023) are appended as rows 72-74.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct F:V for the rows whose paired fixture data was swapped ---
# Row 40: Young Boys 4-1 Lugano
$ws.Range("F40").Value = 'Young Boys'
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 'Lugano'
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = 1.65
$ws.Range("K40").Value = '16/09/2023 17:13'
$ws.Range("L40").Value = 1.67
$ws.Range("M40").Value = '24/09/2023 16:21'
$ws.Range("N40").Value = 4.33
$ws.Range("O40").Value = '16/09/2023 17:13'
$ws.Range("P40").Value = 4.45
$ws.Range("Q40").Value = '24/09/2023 16:29'
$ws.Range("R40").Value = 4.37
$ws.Range("S40").Value = '16/09/2023 17:13'
$ws.Range("T40").Value = 4.68
$ws.Range("U40").Value = '24/09/2023 16:21'
$ws.Range("V40").Value = 'https://www.betexplorer.com/football/switzerland/super-league/young-boys-lugano/hGIXUKka/'

# Row 41: Luzern 2-0 Servette
$ws.Range("F41").Value = 'Luzern'
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 'Servette'
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 2.09
$ws.Range("K41").Value = '16/09/2023 17:13'
$ws.Range("L41").Value = 2.16
$ws.Range("M41").Value = '24/09/2023 16:28'
$ws.Range("N41").Value = 3.66
$ws.Range("O41").Value = '16/09/2023 17:13'
$ws.Range("P41").Value = 3.76
$ws.Range("Q41").Value = '24/09/2023 16:28'
$ws.Range("R41").Value = 3.23
$ws.Range("S41").Value = '16/09/2023 17:13'
$ws.Range("T41").Value = 3.29
$ws.Range("U41").Value = '24/09/2023 16:28'
$ws.Range("V41").Value = 'https://www.betexplorer.com/football/switzerland/super-league/luzern-servette/AiTSV0zg/'

# Row 48: Winterthur 2-3 Lugano
$ws.Range("F48").Value = 'Winterthur'
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 'Lugano'
$ws.Range("I48").Value = 3
$ws.Range("J48").Value = 3.05
$ws.Range("K48").Value = '27/09/2023 19:42'
$ws.Range("L48").Value = 2.58
$ws.Range("M48").Value = '30/09/2023 17:59'
$ws.Range("N48").Value = 3.54
$ws.Range("O48").Value = '27/09/2023 19:42'
$ws.Range("P48").Value = 3.68
$ws.Range("Q48").Value = '30/09/2023 17:59'
$ws.Range("R48").Value = 2.31
$ws.Range("S48").Value = '27/09/2023 19:42'
$ws.Range("T48").Value = 2.68
$ws.Range("U48").Value = '30/09/2023 17:59'
$ws.Range("V48").Value = 'https://www.betexplorer.com/football/switzerland/super-league/winterthur-lugano/6Jr7qK4h/'

# Row 49: Servette 2-1 Lausanne
$ws.Range("F49").Value = 'Servette'
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 'Lausanne'
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = 1.74
$ws.Range("K49").Value = '27/09/2023 19:42'
$ws.Range("L49").Value = 1.88
$ws.Range("M49").Value = '30/09/2023 17:34'
$ws.Range("N49").Value = 3.96
$ws.Range("O49").Value = '27/09/2023 19:42'
$ws.Range("P49").Value = 3.95
$ws.Range("Q49").Value = '30/09/2023 17:54'
$ws.Range("R49").Value = 4.23
$ws.Range("S49").Value = '27/09/2023 19:42'
$ws.Range("T49").Value = 3.97
$ws.Range("U49").Value = '30/09/2023 17:54'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/switzerland/super-league/servette-lausanne/lrl2p0kn/'

# Row 52: Luzern 1-4 Zurich
$ws.Range("F52").Value = 'Luzern'
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 'Zurich'
$ws.Range("I52").Value = 4
$ws.Range("J52").Value = 2.26
$ws.Range("K52").Value = '28/09/2023 19:42'
$ws.Range("L52").Value = 2.42
$ws.Range("M52").Value = '01/10/2023 16:20'
$ws.Range("N52").Value = 3.65
$ws.Range("O52").Value = '28/09/2023 19:42'
$ws.Range("P52").Value = 3.61
$ws.Range("Q52").Value = '01/10/2023 16:20'
$ws.Range("R52").Value = 2.9
$ws.Range("S52").Value = '28/09/2023 19:42'
$ws.Range("T52").Value = 2.92
$ws.Range("U52").Value = '01/10/2023 16:20'
$ws.Range("V52").Value = 'https://www.betexplorer.com/football/switzerland/super-league/luzern-zurich/zazOux4H/'

# Row 53: Basel 0-3 Lausanne Ouchy
$ws.Range("F53").Value = 'Basel'
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 'Lausanne Ouchy'
$ws.Range("I53").Value = 3
$ws.Range("J53").Value = 1.66
$ws.Range("K53").Value = '28/09/2023 19:42'
$ws.Range("L53").Value = 1.78
$ws.Range("M53").Value = '01/10/2023 16:29'
$ws.Range("N53").Value = 4.21
$ws.Range("O53").Value = '28/09/2023 19:42'
$ws.Range("P53").Value = 4.12
$ws.Range("Q53").Value = '01/10/2023 16:29'
$ws.Range("R53").Value = 4.81
$ws.Range("S53").Value = '28/09/2023 19:42'
$ws.Range("T53").Value = 4.3
$ws.Range("U53").Value = '01/10/2023 16:26'
$ws.Range("V53").Value = 'https://www.betexplorer.com/football/switzerland/super-league/basel-lausanne-ouchy/hzoJtIlB/'

# Row 66: Yverdon 1-1 Winterthur
$ws.Range("F66").Value = 'Yverdon'
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 'Winterthur'
$ws.Range("I66").Value = 1
$ws.Range("J66").Value = 2.65
$ws.Range("K66").Value = '22/10/2023 16:42'
$ws.Range("L66").Value = 2.47
$ws.Range("M66").Value = '28/10/2023 17:59'
$ws.Range("N66").Value = 3.64
$ws.Range("O66").Value = '22/10/2023 16:42'
$ws.Range("P66").Value = 3.73
$ws.Range("Q66").Value = '28/10/2023 17:57'
$ws.Range("R66").Value = 2.44
$ws.Range("S66").Value = '22/10/2023 16:42'
$ws.Range("T66").Value = 2.78
$ws.Range("U66").Value = '28/10/2023 17:59'
$ws.Range("V66").Value = 'https://www.betexplorer.com/football/switzerland/super-league/yverdon-winterthur/8QZZ7fmA/'

# Row 67: St. Gallen 3-1 Grasshoppers
$ws.Range("F67").Value = 'St. Gallen'
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 'Grasshoppers'
$ws.Range("I67").Value = 1
$ws.Range("J67").Value = 1.39
$ws.Range("K67").Value = '22/10/2023 16:42'
$ws.Range("L67").Value = 1.62
$ws.Range("M67").Value = '28/10/2023 17:36'
$ws.Range("N67").Value = 5.14
$ws.Range("O67").Value = '22/10/2023 16:42'
$ws.Range("P67").Value = 4.56
$ws.Range("Q67").Value = '28/10/2023 17:59'
$ws.Range("R67").Value = 6.33
$ws.Range("S67").Value = '22/10/2023 16:42'
$ws.Range("T67").Value = 5.02
$ws.Range("U67").Value = '28/10/2023 17:59'
$ws.Range("V67").Value = 'https://www.betexplorer.com/football/switzerland/super-league/st-gallen-grasshoppers/fNVV8zY3/'

# Row 70: Lugano 1-1 Young Boys
$ws.Range("F70").Value = 'Lugano'
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 'Young Boys'
$ws.Range("I70").Value = 1
$ws.Range("J70").Value = 2.95
$ws.Range("K70").Value = '22/10/2023 20:15'
$ws.Range("L70").Value = 2.73
$ws.Range("M70").Value = '29/10/2023 16:28'
$ws.Range("N70").Value = 3.68
$ws.Range("O70").Value = '22/10/2023 20:15'
$ws.Range("P70").Value = 3.66
$ws.Range("Q70").Value = '29/10/2023 16:28'
$ws.Range("R70").Value = 2.22
$ws.Range("S70").Value = '22/10/2023 20:15'
$ws.Range("T70").Value = 2.53
$ws.Range("U70").Value = '29/10/2023 16:29'
$ws.Range("V70").Value = 'https://www.betexplorer.com/football/switzerland/super-league/lugano-young-boys/Y7sNAd3i/'

# Row 71: Servette 4-2 Luzern
$ws.Range("F71").Value = 'Servette'
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 'Luzern'
$ws.Range("I71").Value = 2
$ws.Range("J71").Value = 1.99
$ws.Range("K71").Value = '22/10/2023 20:15'
$ws.Range("L71").Value = 1.93
$ws.Range("M71").Value = '29/10/2023 16:21'
$ws.Range("N71").Value = 3.86
$ws.Range("O71").Value = '22/10/2023 20:15'
$ws.Range("P71").Value = 3.87
$ws.Range("Q71").Value = '29/10/2023 16:21'
$ws.Range("R71").Value = 3.52
$ws.Range("S71").Value = '22/10/2023 20:15'
$ws.Range("T71").Value = 3.9
$ws.Range("U71").Value = '29/10/2023 16:21'
$ws.Range("V71").Value = 'https://www.betexplorer.com/football/switzerland/super-league/servette-luzern/CbtR9GIc/'

# --- Append the 3 new fixture rows, copying row-71 formatting first (bold/bordered index, date format) ---
$ws.Range("A71:V71").Copy()
$ws.Range("A72:V74").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 72: Winterthur 1-4 Young Boys
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = 'switzerland'
$ws.Range("C72").Value = 'super-league'
$ws.Range("D72").Value = '2023-2024'
$ws.Range("E72").Value = 45234.75
$ws.Range("F72").Value = 'Winterthur'
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 'Young Boys'
$ws.Range("I72").Value = 4
$ws.Range("J72").Value = 3.25
$ws.Range("K72").Value = '29/10/2023 16:42'
$ws.Range("L72").Value = 3.47
$ws.Range("M72").Value = '04/11/2023 17:57'
$ws.Range("N72").Value = 4.23
$ws.Range("O72").Value = '29/10/2023 16:42'
$ws.Range("P72").Value = 3.89
$ws.Range("Q72").Value = '04/11/2023 17:57'
$ws.Range("R72").Value = 1.93
$ws.Range("S72").Value = '29/10/2023 16:42'
$ws.Range("T72").Value = 2.04
$ws.Range("U72").Value = '04/11/2023 17:57'
$ws.Range("V72").Value = 'https://www.betexplorer.com/football/switzerland/super-league/winterthur-young-boys/WOvFbjAd/'

# Row 73: Lausanne 3-1 Lugano
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = 'switzerland'
$ws.Range("C73").Value = 'super-league'
$ws.Range("D73").Value = '2023-2024'
$ws.Range("E73").Value = 45234.75
$ws.Range("F73").Value = 'Lausanne'
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 'Lugano'
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = 2.24
$ws.Range("K73").Value = '29/10/2023 16:42'
$ws.Range("L73").Value = 2.34
$ws.Range("M73").Value = '04/11/2023 17:52'
$ws.Range("N73").Value = 3.67
$ws.Range("O73").Value = '29/10/2023 16:42'
$ws.Range("P73").Value = 3.67
$ws.Range("Q73").Value = '04/11/2023 17:52'
$ws.Range("R73").Value = 3.08
$ws.Range("S73").Value = '29/10/2023 16:42'
$ws.Range("T73").Value = 3
$ws.Range("U73").Value = '04/11/2023 17:52'
$ws.Range("V73").Value = 'https://www.betexplorer.com/football/switzerland/super-league/lausanne-lugano/Mkku5hXS/'

# Row 74: Zurich 0-2 Servette
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = 'switzerland'
$ws.Range("C74").Value = 'super-league'
$ws.Range("D74").Value = '2023-2024'
$ws.Range("E74").Value = 45234.85416666666
$ws.Range("F74").Value = 'Zurich'
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 'Servette'
$ws.Range("I74").Value = 2
$ws.Range("J74").Value = 2.08
$ws.Range("K74").Value = '29/10/2023 16:42'
$ws.Range("L74").Value = 2.58
$ws.Range("M74").Value = '04/11/2023 20:21'
$ws.Range("N74").Value = 3.69
$ws.Range("O74").Value = '29/10/2023 16:42'
$ws.Range("P74").Value = 3.46
$ws.Range("Q74").Value = '04/11/2023 20:21'
$ws.Range("R74").Value = 3.41
$ws.Range("S74").Value = '29/10/2023 16:42'
$ws.Range("T74").Value = 2.81
$ws.Range("U74").Value = '04/11/2023 20:21'
$ws.Range("V74").Value = 'https://www.betexplorer.com/football/switzerland/super-league/zurich-servette/YTIcNZnc/'

